$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Page")

# --- Two new rows for the Super Admin credentials (write B10 before A10 so the
#     shared-string table grows in the same order the original edit session used) ---
$ws.Range("B10").Value() = "automate@workstreets.com"
$ws.Range("A10").Value() = "Super Admin Email"
$ws.Range("A11").Value() = "Super Admin Pwd"
$ws.Range("B11").Value() = "admin@123"

# --- Candidate Email (row 8) now points to a different mailbox ---
$ws.Range("B8").Value() = "can06@mailinator.com"

# Carry over the label-column formatting used by the other rows in this block
$ws.Range("A10").Style = $ws.Range("A9").Style
$ws.Range("A11").Style = $ws.Range("A9").Style
$ws.Range("B10").Style = $ws.Range("B9").Style
$ws.Range("B11").Style = $ws.Range("B9").Style

# --- Hyperlinks: rebuild the whole set so ids/order line up with the new rows ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:admin@123")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:wipro@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:adminsiva@nada.email")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:admin@123")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:admin@123")
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:can06@mailinator.com")
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:admin@123")
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:automate@workstreets.com")
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:admin@123")

# Hyperlinks.Add() re-stamps the cell style it touches, so reapply the
# original per-cell formatting once more after rebuilding the links.
$refStyle = $ws.Range("B9").Style
$ws.Range("B2").Style = $ws.Range("B3").Style
$ws.Range("B3").Style = $ws.Range("B2").Style
$ws.Range("B4").Style = $refStyle
$ws.Range("B5").Style = $refStyle
$ws.Range("B7").Style = $refStyle
$ws.Range("B8").Style = $refStyle
$ws.Range("B9").Style = $refStyle
$ws.Range("B10").Style = $refStyle
$ws.Range("B11").Style = $refStyle

# --- Keep the active selection on the row that changed, like the source file ---
$ws.Range("B8").Select()
